{"js": "// Office.js (Word JavaScript API) script.\n// Rewrites the \"Stellar Symphony\" astronomy essay into\n// \"The Art of Governance\" political-science essay, per the commit diff.\n//\n// Strategy: the document body has 7 paragraphs:\n//   0: Title\n//   1: Author name\n//   2: Author email (split across several runs)\n//   3: (blank)\n//   4: Intro / body essay (one paragraph, several <w:br/> line breaks inside)\n//   5: \"Summary\" heading\n//   6: Summary paragraph\n// We rewrite the text paragraph-by-paragraph, preserving each paragraph's\n// existing formatting (font/size/color) by using insertText(..., \"Replace\")\n// on the paragraph itself, and finish by appending a new blank paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Title ---------------------------------------------------------\nparagraphs.items[0].insertText(\n  \"The Art of Governance: Understanding the Dynamics of Power\",\n  Word.InsertLocation.replace\n);\n\n// --- 2. Author name -----------------------------------------------------\nparagraphs.items[1].insertText(\n  \"Tristan Freemantle\",\n  Word.InsertLocation.replace\n);\n\n// --- 3. Author email (replace the whole paragraph's text in one shot,\n//        since it previously was split into five runs: \"alexie\" + \".\" +\n//        \"cunningham@astrotechnologies\" + \".\" + \"org\") -------------------\nparagraphs.items[2].insertText(\n  \"freemantlet@scholarly.edu\",\n  Word.InsertLocation.replace\n);\n\n// --- 4. Intro / body essay paragraph -------------------------------------\n// \"\\v\" is the Office.js text-run encoding for a manual line break (<w:br/>).\nconst introText =\n  \"A nation's governance is akin to navigating a boundless sea of human \" +\n  \"interaction, a complex web of laws, policies, and institutions that \" +\n  \"define the lifeblood of society. Within this dynamic realm, the study \" +\n  \"of governance delves into the intricate mechanisms of power, authority, \" +\n  \"and decision-making that shape our world. To comprehend the essence of \" +\n  \"governance, we must embark on a journey of exploration, examining the \" +\n  \"interplay of power structures, the role of citizens, and the \" +\n  \"challenges of leadership in a globalized world.\\v\\v\" +\n  \"The foundation of governance lies in the distribution and exercise of \" +\n  \"power. Power, in its myriad forms, flows through the veins of society, \" +\n  \"influencing everything from resource allocation to decision-making. \" +\n  \"Unraveling the intricacies of power dynamics involves understanding \" +\n  \"who holds it, how it is wielded, and how it is legitimized. The study \" +\n  \"of governance sheds light on the structures--formal and \" +\n  \"informal--through which power is exercised, ensuring a balance \" +\n  \"between different societal interests\\v\\v\" +\n  \"Simultaneously, governance encompasses the role of citizens in shaping \" +\n  \"their own destiny. Active participation is the lifeblood of a healthy \" +\n  \"democracy, with citizens acting as both subjects and agents of change. \" +\n  \"The study of governance explores the mechanisms through which \" +\n  \"citizens can engage in decision-making processes, hold their leaders \" +\n  \"accountable, and influence policies that impact their lives. \" +\n  \"Understanding the rights, responsibilities, and avenues for civic \" +\n  \"engagement empowers individuals to become active participants in the \" +\n  \"governance of their nation.\\v\\v\" +\n  \"Introduction Continued:\\v\\v\" +\n  \"Leadership, an integral component of governance, demands the ability \" +\n  \"to navigate the complexities of power and the aspirations of citizens. \" +\n  \"Leaders, tasked with the responsibility of guiding their nation's \" +\n  \"course, must possess a keen understanding of the diverse needs and \" +\n  \"interests of their people. The study of governance examines the \" +\n  \"qualities and skills effective leaders need--vision, integrity, \" +\n  \"strategic thinking--to steer their nation through turbulent waters. \" +\n  \"It also explores the challenges they face, from global crises to \" +\n  \"domestic unrest, and the strategies they employ to address these \" +\n  \"challenges.\";\n\nparagraphs.items[4].insertText(introText, Word.InsertLocation.replace);\n\n// --- 5. \"Summary\" heading stays the same; Summary body paragraph --------\nconst summaryText =\n  \"The study of governance is an exploration of the intricate dynamics \" +\n  \"of power, citizenship, and leadership. It delves into the mechanisms \" +\n  \"through which power is distributed and exercised, the role of \" +\n  \"citizens in shaping their own governance, and the qualities and \" +\n  \"challenges of effective leadership. Understanding governance enables \" +\n  \"us to comprehend the functioning of our political systems, appreciate \" +\n  \"the significance of citizen participation, and critically evaluate \" +\n  \"the decisions that shape our world. Ultimately, it empowers us to \" +\n  \"become informed and active citizens, capable of shaping the future \" +\n  \"of our societies.\";\n\nparagraphs.items[6].insertText(summaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// --- 6. Add a new blank paragraph at the very end of the document -------\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Rewrites the \"Stellar Symphony\" astronomy essay into\n# \"The Art of Governance\" political-science essay, per the commit diff.\n#\n# The document has 7 paragraphs:\n#   1: Title\n#   2: Author name\n#   3: Author email (originally split across several runs)\n#   4: (blank)\n#   5: Intro / body essay (one paragraph, several manual line breaks inside)\n#   6: \"Summary\" heading\n#   7: Summary paragraph\n#\n# We rewrite each paragraph's text in place (preserving its formatting) by\n# assigning to a fresh Document.Range(start, end) built from the paragraph's\n# own Range boundaries -- this replaces the full paragraph text even when it\n# originally spanned multiple runs. [char]11 is a manual line break (<w:br/>).\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($paraIndex, $newText) {\n    $p = $d.Paragraphs($paraIndex)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark from the replacement range.\n    $target = $d.Range($r.Start, $r.End - 1)\n    $target.Text = $newText\n}\n\n# --- 1. Title -------------------------------------------------------------\nSet-ParagraphText 1 \"The Art of Governance: Understanding the Dynamics of Power\"\n\n# --- 2. Author name ---------------------------------------------------------\nSet-ParagraphText 2 \"Tristan Freemantle\"\n\n# --- 3. Author email --------------------------------------------------------\nSet-ParagraphText 3 \"freemantlet@scholarly.edu\"\n\n# --- 4. Intro / body essay paragraph ----------------------------------------\n$br = [char]11\n$introText = (\n  \"A nation's governance is akin to navigating a boundless sea of human \" +\n  \"interaction, a complex web of laws, policies, and institutions that \" +\n  \"define the lifeblood of society. Within this dynamic realm, the study \" +\n  \"of governance delves into the intricate mechanisms of power, authority, \" +\n  \"and decision-making that shape our world. To comprehend the essence of \" +\n  \"governance, we must embark on a journey of exploration, examining the \" +\n  \"interplay of power structures, the role of citizens, and the \" +\n  \"challenges of leadership in a globalized world.\" + $br + $br +\n  \"The foundation of governance lies in the distribution and exercise of \" +\n  \"power. Power, in its myriad forms, flows through the veins of society, \" +\n  \"influencing everything from resource allocation to decision-making. \" +\n  \"Unraveling the intricacies of power dynamics involves understanding \" +\n  \"who holds it, how it is wielded, and how it is legitimized. The study \" +\n  \"of governance sheds light on the structures--formal and \" +\n  \"informal--through which power is exercised, ensuring a balance \" +\n  \"between different societal interests\" + $br + $br +\n  \"Simultaneously, governance encompasses the role of citizens in shaping \" +\n  \"their own destiny. Active participation is the lifeblood of a healthy \" +\n  \"democracy, with citizens acting as both subjects and agents of change. \" +\n  \"The study of governance explores the mechanisms through which \" +\n  \"citizens can engage in decision-making processes, hold their leaders \" +\n  \"accountable, and influence policies that impact their lives. \" +\n  \"Understanding the rights, responsibilities, and avenues for civic \" +\n  \"engagement empowers individuals to become active participants in the \" +\n  \"governance of their nation.\" + $br + $br +\n  \"Introduction Continued:\" + $br + $br +\n  \"Leadership, an integral component of governance, demands the ability \" +\n  \"to navigate the complexities of power and the aspirations of citizens. \" +\n  \"Leaders, tasked with the responsibility of guiding their nation's \" +\n  \"course, must possess a keen understanding of the diverse needs and \" +\n  \"interests of their people. The study of governance examines the \" +\n  \"qualities and skills effective leaders need--vision, integrity, \" +\n  \"strategic thinking--to steer their nation through turbulent waters. \" +\n  \"It also explores the challenges they face, from global crises to \" +\n  \"domestic unrest, and the strategies they employ to address these \" +\n  \"challenges.\"\n)\nSet-ParagraphText 5 $introText\n\n# --- 5. \"Summary\" heading stays the same; Summary body paragraph -----------\n$summaryText = (\n  \"The study of governance is an exploration of the intricate dynamics \" +\n  \"of power, citizenship, and leadership. It delves into the mechanisms \" +\n  \"through which power is distributed and exercised, the role of \" +\n  \"citizens in shaping their own governance, and the qualities and \" +\n  \"challenges of effective leadership. Understanding governance enables \" +\n  \"us to comprehend the functioning of our political systems, appreciate \" +\n  \"the significance of citizen participation, and critically evaluate \" +\n  \"the decisions that shape our world. Ultimately, it empowers us to \" +\n  \"become informed and active citizens, capable of shaping the future \" +\n  \"of our societies.\"\n)\nSet-ParagraphText 7 $summaryText\n\n# --- 6. Add a new blank paragraph at the very end of the document ----------\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$lastParagraph.Range.InsertParagraphAfter()\n"}
